$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.023776695207828
$ws.Cells.Item(2, 4).Value = 1.028076240274174
$ws.Cells.Item(2, 5).Value = 1.027401033497069
$ws.Cells.Item(2, 6).Value = 1.022247711018488
$ws.Cells.Item(2, 9).Value = 1.028686702996859
$ws.Cells.Item(2, 10).Value = 1.028955015429087
$ws.Cells.Item(2, 11).Value = 1.03089411017277
$ws.Cells.Item(2, 12).Value = 1.030220870343491
$ws.Cells.Item(2, 13).Value = 1.025082653158573
$ws.Cells.Item(2, 14).Value = 1.013509824634163
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.025362752275332
$ws.Cells.Item(3, 4).Value = 1.029259254878176
$ws.Cells.Item(3, 5).Value = 1.028935580638563
$ws.Cells.Item(3, 6).Value = 1.024480971519728
$ws.Cells.Item(3, 9).Value = 1.029000792948254
$ws.Cells.Item(3, 10).Value = 1.030176401680108
$ws.Cells.Item(3, 11).Value = 1.031883971621403
$ws.Cells.Item(3, 12).Value = 1.031561170960843
$ws.Cells.Item(3, 13).Value = 1.027118645288076
$ws.Cells.Item(3, 14).Value = 1.01393068903879
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.026385986910749
$ws.Cells.Item(4, 4).Value = 1.030021941902366
$ws.Cells.Item(4, 5).Value = 1.029925896153531
$ws.Cells.Item(4, 6).Value = 1.025922337414584
$ws.Cells.Item(4, 9).Value = 1.029201498046556
$ws.Cells.Item(4, 10).Value = 1.030963386428654
$ws.Cells.Item(4, 11).Value = 1.032521119099022
$ws.Cells.Item(4, 12).Value = 1.032425319314288
$ws.Cells.Item(4, 13).Value = 1.028432057511766
$ws.Cells.Item(4, 14).Value = 1.014201464791786
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.026815438853932
$ws.Cells.Item(5, 4).Value = 1.030341914195507
$ws.Cells.Item(5, 5).Value = 1.030341606110571
$ws.Cells.Item(5, 6).Value = 1.026527427069315
$ws.Cells.Item(5, 9).Value = 1.029285271708223
$ws.Cells.Item(5, 10).Value = 1.031293448158413
$ws.Cells.Item(5, 11).Value = 1.032788180804951
$ws.Cells.Item(5, 12).Value = 1.032787873491131
$ws.Cells.Item(5, 13).Value = 1.028983279551106
$ws.Cells.Item(5, 14).Value = 1.014314931232814
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.026887504158777
$ws.Cells.Item(6, 4).Value = 1.030395600447497
$ws.Cells.Item(6, 5).Value = 1.030411369740169
$ws.Cells.Item(6, 6).Value = 1.026628974569171
$ws.Cells.Item(6, 9).Value = 1.029299302427204
$ws.Cells.Item(6, 10).Value = 1.031348821175944
$ws.Cells.Item(6, 11).Value = 1.032832975277508
$ws.Cells.Item(6, 12).Value = 1.032848705252878
$ws.Cells.Item(6, 13).Value = 1.029075777927626
$ws.Cells.Item(6, 14).Value = 1.014333961295409
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.026391728061997
$ws.Cells.Item(7, 4).Value = 1.030026219973189
$ws.Cells.Item(7, 5).Value = 1.029931453303379
$ws.Cells.Item(7, 6).Value = 1.025930425998529
$ws.Cells.Item(7, 9).Value = 1.029202619798024
$ws.Cells.Item(7, 10).Value = 1.030967799804982
$ws.Cells.Item(7, 11).Value = 1.032524690702325
$ws.Cells.Item(7, 12).Value = 1.032430166647642
$ws.Cells.Item(7, 13).Value = 1.028439426615142
$ws.Cells.Item(7, 14).Value = 1.014202982375471
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.024313350900489
$ws.Cells.Item(8, 4).Value = 1.028476631691251
$ws.Cells.Item(8, 5).Value = 1.027920195268045
$ws.Cells.Item(8, 6).Value = 1.023003234276531
$ws.Cells.Item(8, 9).Value = 1.028793377531809
$ws.Cells.Item(8, 10).Value = 1.029368485860404
$ws.Cells.Item(8, 11).Value = 1.031229339493956
$ws.Cells.Item(8, 12).Value = 1.030674483692468
$ws.Cells.Item(8, 13).Value = 1.025771571061255
$ws.Cells.Item(8, 14).Value = 1.013652381302716
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.02062697436999
$ws.Cells.Item(9, 4).Value = 1.025724160925854
$ws.Cells.Item(9, 5).Value = 1.024355296707677
$ws.Cells.Item(9, 6).Value = 1.017815617315894
$ws.Cells.Item(9, 9).Value = 1.028052692111552
$ws.Cells.Item(9, 10).Value = 1.02652423905188
$ws.Cells.Item(9, 11).Value = 1.028920628268163
$ws.Cells.Item(9, 12).Value = 1.027556334871884
$ws.Cells.Item(9, 13).Value = 1.021038678766925
$ws.Cells.Item(9, 14).Value = 1.012670096738093
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.018152334236882
$ws.Cells.Item(10, 4).Value = 1.023873844663236
$ws.Cells.Item(10, 5).Value = 1.021963875212398
$ws.Cells.Item(10, 6).Value = 1.014335705578481
$ws.Cells.Item(10, 9).Value = 1.027545538850753
$ws.Cells.Item(10, 10).Value = 1.024609829127879
$ws.Cells.Item(10, 11).Value = 1.027363335487665
$ws.Cells.Item(10, 12).Value = 1.025460384387471
$ws.Cells.Item(10, 13).Value = 1.017860536854615
$ws.Cells.Item(10, 14).Value = 1.012006893379097
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.017076537603285
$ws.Cells.Item(11, 4).Value = 1.023068861929922
$ws.Cells.Item(11, 5).Value = 1.020924660276648
$ws.Cells.Item(11, 6).Value = 1.012823397272293
$ws.Cells.Item(11, 9).Value = 1.027322716678864
$ws.Cells.Item(11, 10).Value = 1.023776376927169
$ws.Cells.Item(11, 11).Value = 1.026684574788531
$ws.Cells.Item(11, 12).Value = 1.024548565340478
$ws.Cells.Item(11, 13).Value = 1.016478601891089
$ws.Cells.Item(11, 14).Value = 1.011717684724444
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.016676280875228
$ws.Cells.Item(12, 4).Value = 1.022769275036922
$ws.Cells.Item(12, 5).Value = 1.020538075371893
$ws.Cells.Item(12, 6).Value = 1.012260802902105
$ws.Cells.Item(12, 9).Value = 1.027239462249102
$ws.Cells.Item(12, 10).Value = 1.023466105320092
$ws.Cells.Item(12, 11).Value = 1.026431774019943
$ws.Cells.Item(12, 12).Value = 1.024209220258865
$ws.Cells.Item(12, 13).Value = 1.015964392332906
$ws.Cells.Item(12, 14).Value = 1.011609948925069
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.016762167363923
$ws.Cells.Item(13, 4).Value = 1.022833563904492
$ws.Cells.Item(13, 5).Value = 1.020621025393884
$ws.Cells.Item(13, 6).Value = 1.012381520544442
$ws.Cells.Item(13, 9).Value = 1.027257342776716
$ws.Cells.Item(13, 10).Value = 1.023532691103945
$ws.Cells.Item(13, 11).Value = 1.026486031545291
$ws.Cells.Item(13, 12).Value = 1.024282040810072
$ws.Cells.Item(13, 13).Value = 1.01607473311048
$ws.Cells.Item(13, 14).Value = 1.01163307276249
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.017043465783094
$ws.Cells.Item(14, 4).Value = 1.023044109895347
$ws.Cells.Item(14, 5).Value = 1.020892716856389
$ws.Cells.Item(14, 6).Value = 1.01277691068884
$ws.Cells.Item(14, 9).Value = 1.02731584482949
$ws.Cells.Item(14, 10).Value = 1.023750743945503
$ws.Cells.Item(14, 11).Value = 1.02666369211225
$ws.Cells.Item(14, 12).Value = 1.024520528417645
$ws.Cells.Item(14, 13).Value = 1.016436115645429
$ws.Cells.Item(14, 14).Value = 1.011708785618229
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.017216695301557
$ws.Cells.Item(15, 4).Value = 1.023173756949469
$ws.Cells.Item(15, 5).Value = 1.021060038354862
$ws.Cells.Item(15, 6).Value = 1.013020409252274
$ws.Cells.Item(15, 9).Value = 1.027351825031658
$ws.Cells.Item(15, 10).Value = 1.023885001589597
$ws.Cells.Item(15, 11).Value = 1.026773064357672
$ws.Cells.Item(15, 12).Value = 1.024667381270156
$ws.Cells.Item(15, 13).Value = 1.016658655559343
$ws.Cells.Item(15, 14).Value = 1.011755393476236
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.018223640027733
$ws.Cells.Item(16, 4).Value = 1.023927187977966
$ws.Cells.Item(16, 5).Value = 1.022032764858801
$ws.Cells.Item(16, 6).Value = 1.014435954256107
$ws.Cells.Item(16, 9).Value = 1.027560258580779
$ws.Cells.Item(16, 10).Value = 1.02466504658101
$ws.Cells.Item(16, 11).Value = 1.027408288031868
$ws.Cells.Item(16, 12).Value = 1.025520807807639
$ws.Cells.Item(16, 13).Value = 1.017952127094568
$ws.Cells.Item(16, 14).Value = 1.012026043878476
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.018854116553399
$ws.Cells.Item(17, 4).Value = 1.024398773949253
$ws.Cells.Item(17, 5).Value = 1.022641924826689
$ws.Cells.Item(17, 6).Value = 1.015322397761303
$ws.Cells.Item(17, 9).Value = 1.027690137900518
$ws.Cells.Item(17, 10).Value = 1.025153133388854
$ws.Cells.Item(17, 11).Value = 1.027805549663539
$ws.Cells.Item(17, 12).Value = 1.026054988986395
$ws.Cells.Item(17, 13).Value = 1.018761920032639
$ws.Cells.Item(17, 14).Value = 1.012195266964877
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.019221452919525
$ws.Cells.Item(18, 4).Value = 1.024673477828124
$ws.Cells.Item(18, 5).Value = 1.022996880327062
$ws.Cells.Item(18, 6).Value = 1.015838918351605
$ws.Cells.Item(18, 9).Value = 1.027765583784259
$ws.Cells.Item(18, 10).Value = 1.025437392948171
$ws.Cells.Item(18, 11).Value = 1.028036837777103
$ws.Cells.Item(18, 12).Value = 1.026366157960475
$ws.Cells.Item(18, 13).Value = 1.0192337030006
$ws.Cells.Item(18, 14).Value = 1.01229377569063
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.019346636131419
$ws.Cells.Item(19, 4).Value = 1.024767083331571
$ws.Cells.Item(19, 5).Value = 1.023117850881047
$ws.Cells.Item(19, 6).Value = 1.016014950110403
$ws.Cells.Item(19, 9).Value = 1.027791256341327
$ws.Cells.Item(19, 10).Value = 1.025534245059318
$ws.Cells.Item(19, 11).Value = 1.028115628813122
$ws.Cells.Item(19, 12).Value = 1.026472189471204
$ws.Cells.Item(19, 13).Value = 1.019394475447746
$ws.Cells.Item(19, 14).Value = 1.012327331478007
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.018786514917523
$ws.Cells.Item(20, 4).Value = 1.024348214999294
$ws.Cells.Item(20, 5).Value = 1.022576604791662
$ws.Cells.Item(20, 6).Value = 1.015227345492714
$ws.Cells.Item(20, 9).Value = 1.027676235229324
$ws.Cells.Item(20, 10).Value = 1.025100811206189
$ws.Cells.Item(20, 11).Value = 1.027762971589177
$ws.Cells.Item(20, 12).Value = 1.0259977189048
$ws.Cells.Item(20, 13).Value = 1.018675094588026
$ws.Cells.Item(20, 14).Value = 1.012177131261293
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.016960648716593
$ws.Cells.Item(21, 4).Value = 1.02298212546981
$ws.Cells.Item(21, 5).Value = 1.020812726475094
$ws.Cells.Item(21, 6).Value = 1.012660502026747
$ws.Cells.Item(21, 9).Value = 1.027298630949413
$ws.Cells.Item(21, 10).Value = 1.023686551991065
$ws.Cells.Item(21, 11).Value = 1.026611394334261
$ws.Cells.Item(21, 12).Value = 1.024450317961568
$ws.Cells.Item(21, 13).Value = 1.016329722504916
$ws.Cells.Item(21, 14).Value = 1.011686498688233
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.015808838448575
$ws.Cells.Item(22, 4).Value = 1.022119846297559
$ws.Cells.Item(22, 5).Value = 1.019700375767605
$ws.Cells.Item(22, 6).Value = 1.011041658556724
$ws.Cells.Item(22, 9).Value = 1.027058388476454
$ws.Cells.Item(22, 10).Value = 1.022793349728666
$ws.Cells.Item(22, 11).Value = 1.025883418800579
$ws.Cells.Item(22, 12).Value = 1.023473608972092
$ws.Cells.Item(22, 13).Value = 1.014849887243879
$ws.Cells.Item(22, 14).Value = 1.011376217754636
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.016419802545104
$ws.Cells.Item(23, 4).Value = 1.022577279772248
$ws.Cells.Item(23, 5).Value = 1.020290375097155
$ws.Cells.Item(23, 6).Value = 1.01190031969506
$ws.Cells.Item(23, 9).Value = 1.027186015062082
$ws.Cells.Item(23, 10).Value = 1.023267237244101
$ws.Cells.Item(23, 11).Value = 1.026269709063324
$ws.Cells.Item(23, 12).Value = 1.023991745846818
$ws.Cells.Item(23, 13).Value = 1.015634879315715
$ws.Cells.Item(23, 14).Value = 1.011540875845924
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.018817062453231
$ws.Cells.Item(24, 4).Value = 1.024371061537861
$ws.Cells.Item(24, 5).Value = 1.022606121207317
$ws.Cells.Item(24, 6).Value = 1.015270297148547
$ws.Cells.Item(24, 9).Value = 1.027682518207468
$ws.Cells.Item(24, 10).Value = 1.025124454686732
$ws.Cells.Item(24, 11).Value = 1.027782212110188
$ws.Cells.Item(24, 12).Value = 1.026023598055744
$ws.Cells.Item(24, 13).Value = 1.018714328984833
$ws.Cells.Item(24, 14).Value = 1.01218532661141
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.021582929670007
$ws.Cells.Item(25, 4).Value = 1.026438398203029
$ws.Cells.Item(25, 5).Value = 1.025279459176653
$ws.Cells.Item(25, 6).Value = 1.019160408335618
$ws.Cells.Item(25, 9).Value = 1.028246516048388
$ws.Cells.Item(25, 10).Value = 1.027262706805107
$ws.Cells.Item(25, 11).Value = 1.029520641937565
$ws.Cells.Item(25, 12).Value = 1.028365422786866
$ws.Cells.Item(25, 13).Value = 1.022266166063971
$ws.Cells.Item(25, 14).Value = 1.012925494122837
